$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets("ProductLoanInput")
$wsOutput = $wb.Worksheets("ProductLoanOutput")

# Update the product name text (was "...UPFRONT", now "...UP-1st")
$wsInput.Range("B1").Value = "2480-RBI-EI-DB-SAR-REC-NOCOM-RNI-CTPD-SAR-MD-TR-2-DATE-VAR-INST-UP-1st"

# Update the short name from numeric 2480 to text "248d"
$wsInput.Range("B2").Value = "248d"

# Mirror the product name change on the output sheet
$wsOutput.Range("B1").Value = "2480-RBI-EI-DB-SAR-REC-NOCOM-RNI-CTPD-SAR-MD-TR-2-DATE-VAR-INST-UP-1st"

# Move the selection on the input sheet to B1 (was A36)
$null = $wsInput.Range("B1").Select()

# Switch the active tab from ProductLoanInput to ProductLoanOutput
$wsOutput.Activate()
